$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("A2:D32")
$dst = $ws.Range("A33")
$src.Copy($dst)
